$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: Column A = id (number), Column B = code (string)
# This reflects the final state of the sheet after the cleaning update,
# which now keeps additional rows rather than removing them.
$data = @(
    @(20, "C931253116052"),
    @(5,  "M931100509004"),
    @(21, "X931412020027"),
    @(6,  "U931252114001"),
    @(7,  "C931100609010"),
    @(13, "G931321113006"),
    @(15, "Y931321309029"),
    @(16, "P931101109055"),
    @(2,  "R931258916027"),
    @(11, "L931412020030"),
    @(19, "R931253116053"),
    @(26, "K931383410019")
)

# Use the existing formatted cell (A2) as the formatting template for
# column A so that newly added rows match the look of the existing ones.
$templateA = $ws.Cells.Item(2, 1)

$row = 2
foreach ($item in $data) {
    $cellA = $ws.Cells.Item($row, 1)
    if ($row -gt 5) {
        $templateA.Copy() | Out-Null
        $cellA.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    }
    $cellA.Value = $item[0]

    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
